# Update sector correlation data for 2016 intra-sector correlations sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 'Multi-Utilities(18)'
$ws.Cells.Item(2, 2).Value = 0.6011530979219659
$ws.Cells.Item(3, 1).Value = 'Road & Rail(22)'
$ws.Cells.Item(3, 2).Value = 0.555478110834445
$ws.Cells.Item(4, 1).Value = 'Air Freight & Logistics(11)'
$ws.Cells.Item(4, 2).Value = 0.5120315985335757
$ws.Cells.Item(5, 1).Value = 'Electric Utilities(28)'
$ws.Cells.Item(5, 2).Value = 0.5111832501130005
$ws.Cells.Item(6, 1).Value = 'Banks(246)'
$ws.Cells.Item(6, 2).Value = 0.4914866227540425
$ws.Cells.Item(7, 1).Value = 'Building Products(23)'
$ws.Cells.Item(7, 2).Value = 0.4817696952719114
$ws.Cells.Item(8, 1).Value = 'Energy Equipment & Services(32)'
$ws.Cells.Item(8, 2).Value = 0.4718498402254432
$ws.Cells.Item(9, 1).Value = 'Metals & Mining(89)'
$ws.Cells.Item(9, 2).Value = 0.4316523242605258
$ws.Cells.Item(10, 1).Value = 'Machinery(85)'
$ws.Cells.Item(10, 2).Value = 0.4152030506263252
$ws.Cells.Item(11, 1).Value = 'Marine(15)'
$ws.Cells.Item(11, 2).Value = 0.3878966859677347
$ws.Cells.Item(12, 1).Value = 'Gas Utilities(12)'
$ws.Cells.Item(12, 2).Value = 0.3704453840804315
$ws.Cells.Item(13, 1).Value = 'Trading Companies & Distributors(25)'
$ws.Cells.Item(13, 2).Value = 0.3654504162682767
$ws.Cells.Item(14, 1).Value = 'Thrifts & Mortgage Finance(47)'
$ws.Cells.Item(14, 2).Value = 0.3547728109197359
$ws.Cells.Item(15, 1).Value = 'Water Utilities(12)'
$ws.Cells.Item(15, 2).Value = 0.3520959880163748
$ws.Cells.Item(16, 1).Value = 'Insurance(75)'
$ws.Cells.Item(16, 2).Value = 0.3180148923140174
$ws.Cells.Item(17, 1).Value = 'Auto Components(21)'
$ws.Cells.Item(17, 2).Value = 0.3055511741165613
$ws.Cells.Item(18, 1).Value = 'Construction & Engineering(20)'
$ws.Cells.Item(18, 2).Value = 0.3023106614107945
$ws.Cells.Item(19, 1).Value = 'Chemicals(51)'
$ws.Cells.Item(19, 2).Value = 0.2921011422361156
$ws.Cells.Item(20, 1).Value = 'Life Sciences Tools & Services(19)'
$ws.Cells.Item(20, 2).Value = 0.2631741957874983
$ws.Cells.Item(21, 1).Value = 'Specialty Retail(58)'
$ws.Cells.Item(21, 2).Value = 0.2409116727380363
$ws.Cells.Item(22, 1).Value = 'Capital Markets(75)'
$ws.Cells.Item(22, 2).Value = 0.2380249774285201
$ws.Cells.Item(23, 1).Value = 'Semiconductors & Semiconductor Equipment(68)'
$ws.Cells.Item(23, 2).Value = 0.2154474925683811
$ws.Cells.Item(24, 1).Value = 'Electrical Equipment(28)'
$ws.Cells.Item(24, 2).Value = 0.1998757633264295
$ws.Cells.Item(25, 1).Value = 'Commercial Services & Supplies(52)'
$ws.Cells.Item(25, 2).Value = 0.1914991355344332
$ws.Cells.Item(26, 1).Value = 'Professional Services(35)'
$ws.Cells.Item(26, 2).Value = 0.1896722757203425
$ws.Cells.Item(27, 1).Value = 'Aerospace & Defense(37)'
$ws.Cells.Item(27, 2).Value = 0.1828968893596951
$ws.Cells.Item(28, 1).Value = 'Hotels, Restaurants & Leisure(50)'
$ws.Cells.Item(28, 2).Value = 0.1787714773155454
$ws.Cells.Item(29, 1).Value = 'Oil, Gas & Consumable Fuels(122)'
$ws.Cells.Item(29, 2).Value = 0.1783105191961414
$ws.Cells.Item(30, 1).Value = 'Pharmaceuticals(48)'
$ws.Cells.Item(30, 2).Value = 0.1745937056713801
$ws.Cells.Item(31, 1).Value = 'Health Care Providers & Services(46)'
$ws.Cells.Item(31, 2).Value = 0.1591865300992264
$ws.Cells.Item(32, 1).Value = 'Communications Equipment(45)'
$ws.Cells.Item(32, 2).Value = 0.1519806083713743
$ws.Cells.Item(33, 1).Value = 'Media(42)'
$ws.Cells.Item(33, 2).Value = 0.1391667115132491
$ws.Cells.Item(34, 1).Value = 'Household Durables(39)'
$ws.Cells.Item(34, 2).Value = 0.1349024315186565
$ws.Cells.Item(35, 1).Value = 'Health Care Equipment & Supplies(83)'
$ws.Cells.Item(35, 2).Value = 0.1174498816020266
$ws.Cells.Item(36, 1).Value = 'Biotechnology(126)'
$ws.Cells.Item(36, 2).Value = 0.1129795663768614
$ws.Cells.Item(37, 1).Value = 'Food Products(44)'
$ws.Cells.Item(37, 2).Value = 0.1015430993239691
$ws.Cells.Item(38, 1).Value = 'IT Services(52)'
$ws.Cells.Item(38, 2).Value = 0.09920933229872297
$ws.Cells.Item(39, 1).Value = 'Software(66)'
$ws.Cells.Item(39, 2).Value = 0.08997060194791422

# Remove now-unused trailing rows (sheet shrank from 41 to 39 data rows)
$ws.Range("A40:B41").Delete()
